$wb = $excel.ActiveWorkbook

# ALC
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 53522.656
$ws.Range("J17").Value = 53522.656
$ws.Range("L17").Value = 160567.968
$ws.Range("N17").Value = -160903.968
$ws.Range("H19").Value = 4221.3335
$ws.Range("I19").Value = 4000.5
$ws.Range("J19").Value = 4284.4287
$ws.Range("K19").Value = 4000.5
$ws.Range("L19").Value = 4284.4287
$ws.Range("M19").Value = -3825.5
$ws.Range("N19").Value = -4634.4287
$ws.Range("H106").Value = 3301.1
$ws.Range("I106").Value = 4169.1665
$ws.Range("K106").Value = 4169.1665
$ws.Range("M106").Value = -3538.1665
$ws.Range("H128").Value = 87392.086
$ws.Range("J128").Value = 87392.086
$ws.Range("L128").Value = 87392.086
$ws.Range("N128").Value = -97352.086
$ws.Range("H131").Value = 2689.9614
$ws.Range("I131").Value = 1953.05
$ws.Range("K131").Value = 5859.15
$ws.Range("M131").Value = -819.1499999999996
$ws.Range("H133").Value = 49999.145
$ws.Range("J133").Value = 49999.145
$ws.Range("L133").Value = 49999.145
$ws.Range("N133").Value = -60119.145
$ws.Range("H137").Value = 3546.1667
$ws.Range("I137").Value = 1811
$ws.Range("J137").Value = 7760.143
$ws.Range("K137").Value = 5433
$ws.Range("L137").Value = 23280.429
$ws.Range("M137").Value = -2883
$ws.Range("N137").Value = -28380.429
$ws.Range("H138").Value = 3675.4243
$ws.Range("I138").Value = 1986.5238
$ws.Range("J138").Value = 4463.5776
$ws.Range("K138").Value = 5959.5714
$ws.Range("L138").Value = 13390.7328
$ws.Range("M138").Value = -819.5713999999998
$ws.Range("N138").Value = -23670.7328

# ARM
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2050.6086
$ws.Range("I32").Value = 1601.3064
$ws.Range("K32").Value = 1601.3064
$ws.Range("M32").Value = -1314.3064
$ws.Range("H61").Value = 3774
$ws.Range("I61").Value = 2936.375
$ws.Range("K61").Value = 2936.375
$ws.Range("M61").Value = -2724.375
$ws.Range("H63").Value = 3624.4167
$ws.Range("I63").Value = 3499.3635
$ws.Range("J63").Value = 5000
$ws.Range("K63").Value = 3499.3635
$ws.Range("L63").Value = 5000
$ws.Range("M63").Value = -2813.3635
$ws.Range("N63").Value = -6372
$ws.Range("H66").Value = 3624.4167
$ws.Range("I66").Value = 3499.3635
$ws.Range("J66").Value = 5000
$ws.Range("K66").Value = 17496.8175
$ws.Range("L66").Value = 25000
$ws.Range("M66").Value = -14064.8175
$ws.Range("N66").Value = -31864
$ws.Range("H88").Value = 2714.5715
$ws.Range("I88").Value = 0
$ws.Range("K88").Value = 0
$ws.Range("M88").ClearContents()
$ws.Range("H91").Value = 2714.5715
$ws.Range("I91").Value = 0
$ws.Range("K91").Value = 0
$ws.Range("M91").ClearContents()
$ws.Range("H136").Value = 3774
$ws.Range("I136").Value = 2936.375
$ws.Range("K136").Value = 8809.125
$ws.Range("M136").Value = -6259.125

# BSM
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H64").Value = 1000.8
$ws.Range("I64").Value = 990
$ws.Range("J64").Value = 1006.2
$ws.Range("K64").Value = 990
$ws.Range("L64").Value = 1006.2
$ws.Range("M64").Value = -765
$ws.Range("N64").Value = -1456.2
$ws.Range("H67").Value = 1000.8
$ws.Range("I67").Value = 990
$ws.Range("J67").Value = 1006.2
$ws.Range("K67").Value = 990
$ws.Range("L67").Value = 1006.2
$ws.Range("M67").Value = -210
$ws.Range("N67").Value = -2566.2
$ws.Range("H99").Value = 2035.2609
$ws.Range("I99").Value = 1640.6
$ws.Range("K99").Value = 1640.6
$ws.Range("M99").Value = -142.5999999999999
$ws.Range("H105").Value = 27401.238
$ws.Range("I105").Value = 60817.883
$ws.Range("J105").Value = 4677.92
$ws.Range("K105").Value = 60817.883
$ws.Range("L105").Value = 4677.92
$ws.Range("M105").Value = -59070.883
$ws.Range("N105").Value = -8171.92

# CRP
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H43").Value = 23328.5
$ws.Range("J43").Value = 23328.5
$ws.Range("L43").Value = 23328.5
$ws.Range("N43").Value = -23696.5
$ws.Range("H55").Value = 16488.25
$ws.Range("I55").Value = 14291
$ws.Range("K55").Value = 14291
$ws.Range("M55").Value = -13976
$ws.Range("H101").Value = 23328.5
$ws.Range("J101").Value = 23328.5
$ws.Range("L101").Value = 23328.5
$ws.Range("N101").Value = -29818.5

# CUL
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 1518.6
$ws.Range("I68").Value = 1265
$ws.Range("K68").Value = 3795
$ws.Range("M68").Value = -2984
$ws.Range("H71").Value = 1518.6
$ws.Range("I71").Value = 1265
$ws.Range("K71").Value = 11385
$ws.Range("M71").Value = -7329
$ws.Range("H92").Value = 989
$ws.Range("I92").Value = 640.4
$ws.Range("J92").Value = 1424.75
$ws.Range("K92").Value = 1921.2
$ws.Range("L92").Value = 4274.25
$ws.Range("M92").Value = -673.1999999999998
$ws.Range("N92").Value = -6770.25
$ws.Range("H106").Value = 5610.4
$ws.Range("J106").Value = 10029
$ws.Range("L106").Value = 30087
$ws.Range("N106").Value = -31979
$ws.Range("H108").Value = 5848.1
$ws.Range("I108").Value = 5622.625
$ws.Range("K108").Value = 16867.875
$ws.Range("M108").Value = -13987.875
$ws.Range("H113").Value = 4631229
$ws.Range("J113").Value = 1999.3334
$ws.Range("L113").Value = 5998.0002
$ws.Range("N113").Value = -10338.0002
$ws.Range("H122").Value = 250751.25
$ws.Range("I122").Value = 1000
$ws.Range("K122").Value = 9000
$ws.Range("M122").Value = -6550

# GSM
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H93").Value = 43999
$ws.Range("J93").Value = 43999
$ws.Range("L93").Value = 43999
$ws.Range("N93").Value = -47743
$ws.Range("H99").Value = 3626.3333
$ws.Range("I99").Value = 3626.3333
$ws.Range("K99").Value = 3626.3333
$ws.Range("M99").Value = -1380.3333
$ws.Range("H113").Value = 393040.8
$ws.Range("I113").Value = 463380.6
$ws.Range("K113").Value = 463380.6
$ws.Range("M113").Value = -461210.6
$ws.Range("H132").Value = 1254750.5
$ws.Range("I132").Value = 3337003
$ws.Range("K132").Value = 10011009
$ws.Range("M132").Value = -10008479

# LTW
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 5201.1113
$ws.Range("I16").Value = 5045.857
$ws.Range("J16").Value = 5744.5
$ws.Range("K16").Value = 5045.857
$ws.Range("L16").Value = 5744.5
$ws.Range("M16").Value = -4875.857
$ws.Range("N16").Value = -6084.5
$ws.Range("H46").Value = 3828.64
$ws.Range("I46").Value = 3218.2354
$ws.Range("K46").Value = 3218.2354
$ws.Range("M46").Value = -3030.2354
$ws.Range("H132").Value = 5129
$ws.Range("I132").Value = 3711.6667
$ws.Range("K132").Value = 11135.0001
$ws.Range("M132").Value = -8605.000100000001

# WVR
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 18408.572
$ws.Range("I81").Value = 2361
$ws.Range("J81").Value = 24827.6
$ws.Range("K81").Value = 4722
$ws.Range("L81").Value = 49655.2
$ws.Range("M81").Value = -3661
$ws.Range("N81").Value = -51777.2
$ws.Range("H84").Value = 18408.572
$ws.Range("I84").Value = 2361
$ws.Range("J84").Value = 24827.6
$ws.Range("K84").Value = 23610
$ws.Range("L84").Value = 248276
$ws.Range("M84").Value = -18306
$ws.Range("N84").Value = -258884
$ws.Range("H122").Value = 31254180
$ws.Range("I122").Value = 50002644
$ws.Range("K122").Value = 150007932
$ws.Range("M122").Value = -150005482
